$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 299
$ws1.Range("F4").Value = 879
$ws1.Range("F5").Value = 44
$ws1.Range("F6").Value = 347
$ws1.Range("F7").Value = 10804
$ws1.Range("F8").Value = 357
$ws1.Range("F12").Value = 139
$ws1.Range("F13").Value = 140
$ws1.Range("F14").Value = 15
$ws1.Range("F16").Value = 40
$ws1.Range("F20").Value = 1082
$ws1.Range("F21").Value = 52
$ws1.Range("F22").Value = 106

# Sheet "全部类型" (all types) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 299
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 879
$ws4.Range("F5").Value = 44
$ws4.Range("F6").Value = 347
$ws4.Range("F7").Value = 10805
$ws4.Range("F8").Value = 357
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 0
$ws4.Range("F12").Value = 0
$ws4.Range("F13").Value = 140
$ws4.Range("F14").Value = 15
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 40
$ws4.Range("F17").Value = 0
$ws4.Range("F19").Value = 0
$ws4.Range("F20").Value = 1082
$ws4.Range("F21").Value = 52
$ws4.Range("F22").Value = 106

$wb.Save()
